$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: numero_carte=1 (number), code_client=C00005 (text), taux_reduction=5 (number)
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "C00005"
$ws.Cells.Item(2, 3).Value = 5

# Row 3: numero_carte=000002 (text, preserve leading zeros), code_client=C00006 (text), taux_reduction=5 (number)
# Format the cell as Text first so the numeric-looking string "000002" is
# stored as a literal string instead of being coerced to the number 2,
# then clear the formatting so the cell keeps the workbook's default style.
$ws.Cells.Item(3, 1).NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = "000002"
$ws.Cells.Item(3, 1).ClearFormats()
$ws.Cells.Item(3, 2).Value = "C00006"
$ws.Cells.Item(3, 3).Value = 5
